# add txt preprocessing step
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 6 (2023-01-15_15-15-30.wav): normalize mechanism casing and
# rewrite injuries / activation_page text after txt preprocessing
$ws.Range("H6").Value = "fall from height"
$ws.Range("I6").Value = "Fracture, Internal Bleeding"
$ws.Range("J6").Value = "Trauma activation for 32-year-old F with fall from height. Injuries: Fracture, Internal Bleeding. Trauma team activation required."

# Row 14 (2023-01-15_14-30-25.wav): same preprocessing applied
$ws.Range("H14").Value = "motor vehicle accident"
$ws.Range("I14").Value = "Head Injury, Chest Trauma, Fracture, Rib Fracture, Skull Fracture"
$ws.Range("J14").Value = "Trauma activation for 25-year-old M with motor vehicle accident. Injuries: Head Injury, Chest Trauma, Fracture, Rib Fracture, Skull Fracture. Trauma team activation required."
